## Automatische test-sync: 2025-06-22 18:53:50
## Adds the new "Aanmelding nieuwsbrief" mail-log entry (row 23) to the
## Logs sheet, extends the two conditionalFormatting ranges to cover it,
## and refreshes the Dashboard category-count table (Afmelding /
## Nieuwsbrief count 1 -> 2, re-sorted by count desc).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new incoming-mail row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A23").Value = "Aanmelding nieuwsbrief"
$logs.Range("B23").Value = "mailmind.test@zohomail.eu"
$logs.Range("C23").Value = "Ik wil me graag inschrijven voor de nieuwsbrief."
$logs.Range("D23").Value = "Afmelding / Nieuwsbrief"
# (no E23 - this mail hasn't been answered yet, same as row 22)
$logs.Range("F23").Value = "2025-06-22 18:53:13"
$logs.Range("G23").Value = "Nee"

# ---------------------------------------------------------------------
# 2) Extend the conditional-formatting ranges from row 22 to row 23
# ---------------------------------------------------------------------
$catRules = $logs.Range("D2:D22").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D23"))
}

$answeredRules = $logs.Range("G2:G22").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G23"))
}

# ---------------------------------------------------------------------
# 3) Dashboard sheet: "Afmelding / Nieuwsbrief" now has 2 occurrences,
#    so it moves up in the count-sorted table; shuffle rows 7-12.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B7").Value = 2

$dash.Range("A8").Value = "Openingstijden / Locatie"
$dash.Range("B8").Value = 1

$dash.Range("A10").Value = "Uitnodiging / Evenement"
$dash.Range("B10").Value = 1

$dash.Range("A12").Value = "Offerte / Prijsaanvraag"
$dash.Range("B12").Value = 1
